$wb = $excel.ActiveWorkbook

# --- Step 1: insert "Đơn thu nợ" sheet + add new "Lương" sheet, matching target tab order ---
$sheetDonThuNo = $wb.Worksheets.Item(3)
$sheetDonThuNo.Name = 'Đơn thu nợ'
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetLuong = $wb.Worksheets.Add($null, $lastSheet)
$sheetLuong.Name = 'Lương'

# --- Step 2: "Đơn sale chính" sheet edits ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 7).Value = 'Nâng mũi'
$ws1.Cells.Item(2, 15).Value = 3000000
$ws1.Cells.Item(2, 16).Value = 17730000
$ws1.Cells.Item(2, 17).Value = 0
$ws1.Cells.Item(2, 24).Value = 0.1
$ws1.Cells.Item(2, 26).Value = 1773000

$ws1.Cells.Item(3, 7).Value = 'Tiểu phẫu'

$ws1.Cells.Item(4, 7).Value = 'Tiêm'

# new row 5
$ws1.Cells.Item(5, 1).Value = 'HD-LUXURY'
$ws1.Cells.Item(5, 2).Value = 545
$ws1.Cells.Item(5, 3).Value = '07-12-2024'
$ws1.Cells.Item(5, 4).Value = 'CẦN THƠ'
$ws1.Cells.Item(5, 5).Value = 'Phạm Thị Tuyết Nhi'
$ws1.Cells.Item(5, 6).Value = 'Khách cũ'
$ws1.Cells.Item(5, 7).Value = 'Tiêm'
$ws1.Cells.Item(5, 8).Value = 'Tiêm Filler'
$ws1.Cells.Item(5, 9).Value = 'Đỗ Thị Huyền Trân'
$ws1.Cells.Item(5, 10).Value = 2800000
$ws1.Cells.Item(5, 11).Value = $null
$ws1.Cells.Item(5, 12).Value = $null
$ws1.Cells.Item(5, 13).Value = 2800000
$ws1.Cells.Item(5, 14).Value = 2800000
$ws1.Cells.Item(5, 15).Value = 0
$ws1.Cells.Item(5, 16).Value = 2800000
$ws1.Cells.Item(5, 17).Value = 0
$ws1.Cells.Item(5, 18).Value = 'Nguyễn Hoàng Yến Quyên'
$ws1.Cells.Item(5, 19).Value = $null
$ws1.Cells.Item(5, 20).Value = $null
$ws1.Cells.Item(5, 21).Value = $null
$ws1.Cells.Item(5, 22).Value = 50000
$ws1.Cells.Item(5, 23).Value = 0
$ws1.Cells.Item(5, 24).Value = 0.1
$ws1.Cells.Item(5, 25).Value = 0
$ws1.Cells.Item(5, 26).Value = 280000
$ws1.Cells.Item(5, 27).Value = 0

# new row 6
$ws1.Cells.Item(6, 1).Value = 'HD-LUXURY'
$ws1.Cells.Item(6, 2).Value = 556
$ws1.Cells.Item(6, 3).Value = '07-16-2024'
$ws1.Cells.Item(6, 4).Value = 'CẦN THƠ'
$ws1.Cells.Item(6, 5).Value = 'Nguyễn Bảo Hân'
$ws1.Cells.Item(6, 6).Value = 'Khách cũ'
$ws1.Cells.Item(6, 7).Value = 'Tiêm'
$ws1.Cells.Item(6, 8).Value = 'Tiêm botox'
$ws1.Cells.Item(6, 9).Value = 'Đỗ Thị Huyền Trân'
$ws1.Cells.Item(6, 10).Value = 1300000
$ws1.Cells.Item(6, 11).Value = $null
$ws1.Cells.Item(6, 12).Value = $null
$ws1.Cells.Item(6, 13).Value = 1300000
$ws1.Cells.Item(6, 14).Value = 1300000
$ws1.Cells.Item(6, 15).Value = 0
$ws1.Cells.Item(6, 16).Value = 1300000
$ws1.Cells.Item(6, 17).Value = 0
$ws1.Cells.Item(6, 18).Value = 'Nguyễn Hoàng Yến Quyên'
$ws1.Cells.Item(6, 19).Value = $null
$ws1.Cells.Item(6, 20).Value = $null
$ws1.Cells.Item(6, 21).Value = $null
$ws1.Cells.Item(6, 22).Value = 0
$ws1.Cells.Item(6, 23).Value = 0
$ws1.Cells.Item(6, 24).Value = 0.1
$ws1.Cells.Item(6, 25).Value = 0
$ws1.Cells.Item(6, 26).Value = 130000
$ws1.Cells.Item(6, 27).Value = 0

# new row 7
$ws1.Cells.Item(7, 1).Value = 'HD-LUXURY'
$ws1.Cells.Item(7, 2).Value = 564
$ws1.Cells.Item(7, 3).Value = '07-17-2024'
$ws1.Cells.Item(7, 4).Value = 'CẦN THƠ'
$ws1.Cells.Item(7, 5).Value = 'Đặng Thị Ngọc Huyền'
$ws1.Cells.Item(7, 6).Value = 'Cá nhân'
$ws1.Cells.Item(7, 7).Value = 'Tiêm'
$ws1.Cells.Item(7, 8).Value = 'Tiêm Filler'
$ws1.Cells.Item(7, 9).Value = 'Đỗ Thị Huyền Trân'
$ws1.Cells.Item(7, 10).Value = 1050000
$ws1.Cells.Item(7, 11).Value = $null
$ws1.Cells.Item(7, 12).Value = $null
$ws1.Cells.Item(7, 13).Value = 1050000
$ws1.Cells.Item(7, 14).Value = 1050000
$ws1.Cells.Item(7, 15).Value = 0
$ws1.Cells.Item(7, 16).Value = 1050000
$ws1.Cells.Item(7, 17).Value = 0
$ws1.Cells.Item(7, 18).Value = 'Nguyễn Hoàng Yến Quyên'
$ws1.Cells.Item(7, 19).Value = $null
$ws1.Cells.Item(7, 20).Value = $null
$ws1.Cells.Item(7, 21).Value = $null
$ws1.Cells.Item(7, 22).Value = 50000
$ws1.Cells.Item(7, 23).Value = 0
$ws1.Cells.Item(7, 24).Value = 0.1
$ws1.Cells.Item(7, 25).Value = 0
$ws1.Cells.Item(7, 26).Value = 105000
$ws1.Cells.Item(7, 27).Value = 0

# new row 8
$ws1.Cells.Item(8, 1).Value = 'HD-LUXURY'
$ws1.Cells.Item(8, 2).Value = 565
$ws1.Cells.Item(8, 3).Value = '07-17-2024'
$ws1.Cells.Item(8, 4).Value = 'CẦN THƠ'
$ws1.Cells.Item(8, 5).Value = 'Sữ Thị Thanh Tuyền'
$ws1.Cells.Item(8, 6).Value = 'Cá nhân'
$ws1.Cells.Item(8, 7).Value = 'Nâng mũi'
$ws1.Cells.Item(8, 8).Value = 'Nâng mũi'
$ws1.Cells.Item(8, 9).Value = 'Đỗ Thị Huyền Trân'
$ws1.Cells.Item(8, 10).Value = 18000000
$ws1.Cells.Item(8, 11).Value = $null
$ws1.Cells.Item(8, 12).Value = $null
$ws1.Cells.Item(8, 13).Value = 18000000
$ws1.Cells.Item(8, 14).Value = 18000000
$ws1.Cells.Item(8, 15).Value = 0
$ws1.Cells.Item(8, 16).Value = 18000000
$ws1.Cells.Item(8, 17).Value = 0
$ws1.Cells.Item(8, 18).Value = 'Lâm Thị Mỹ Hằng'
$ws1.Cells.Item(8, 19).Value = 'Phạm Thanh Hoàng'
$ws1.Cells.Item(8, 20).Value = 'Lâm Hoàng Phú'
$ws1.Cells.Item(8, 21).Value = $null
$ws1.Cells.Item(8, 22).Value = 100000
$ws1.Cells.Item(8, 23).Value = 50000
$ws1.Cells.Item(8, 24).Value = 0.1
$ws1.Cells.Item(8, 25).Value = 0
$ws1.Cells.Item(8, 26).Value = 1800000
$ws1.Cells.Item(8, 27).Value = 0

# new Tổng row 9 (was row 5)
$ws1.Cells.Item(9, 1).Value = 'Tổng'
$ws1.Cells.Item(9, 2).Value = 7
$ws1.Cells.Item(9, 10).Value = 45780000
$ws1.Cells.Item(9, 12).Value = 0
$ws1.Cells.Item(9, 13).Value = 45780000
$ws1.Cells.Item(9, 14).Value = 42780000
$ws1.Cells.Item(9, 15).Value = 3000000
$ws1.Cells.Item(9, 16).Value = 45780000
$ws1.Cells.Item(9, 17).Value = 0
$ws1.Cells.Item(9, 22).Value = 400000
$ws1.Cells.Item(9, 23).Value = 100000
$ws1.Cells.Item(9, 24).Value = 0.64
$ws1.Cells.Item(9, 25).Value = 0
$ws1.Cells.Item(9, 26).Value = 4356000
$ws1.Cells.Item(9, 27).Value = 0

# --- Step 3: "Đơn sale phụ" sheet edits ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2, 7).Value = 'Tiểu phẫu'
$ws2.Cells.Item(2, 24).Value = 0.13
$ws2.Cells.Item(2, 26).Value = 700000

$ws2.Cells.Item(3, 7).Value = 'Tiêm'

$ws2.Cells.Item(4, 7).Value = 'Phun xăm'
$ws2.Cells.Item(4, 15).Value = 500000
$ws2.Cells.Item(4, 16).Value = 1500000
$ws2.Cells.Item(4, 17).Value = 0
$ws2.Cells.Item(4, 27).Value = 30000

$ws2.Cells.Item(5, 7).Value = 'Phun xăm'

$ws2.Cells.Item(6, 7).Value = 'Nâng mũi'
$ws2.Cells.Item(6, 24).Value = 0.1
$ws2.Cells.Item(6, 26).Value = 950000

$ws2.Cells.Item(7, 15).Value = 500000
$ws2.Cells.Item(7, 16).Value = 20900000
$ws2.Cells.Item(7, 17).Value = 8500000
$ws2.Cells.Item(7, 24).Value = 0.33
$ws2.Cells.Item(7, 26).Value = 1830000
$ws2.Cells.Item(7, 27).Value = 248000

# --- Step 4: "Đơn thu nợ" sheet (repurposed sheet) content ---
$ws3 = $sheetDonThuNo

# row 1 (headers)
$ws3.Cells.Item(1, 1).Value = 'Ngày thực hiện'
$ws3.Cells.Item(1, 2).Value = 'Ngày thu'
$ws3.Cells.Item(1, 3).Value = 'notion id'
$ws3.Cells.Item(1, 4).Value = 'Tiền tố'
$ws3.Cells.Item(1, 5).Value = 'Mã đơn thu nợ'
$ws3.Cells.Item(1, 6).Value = 'Cơ sở'
$ws3.Cells.Item(1, 7).Value = 'id đơn nợ'
$ws3.Cells.Item(1, 8).Value = 'Lượng thu'
$ws3.Cells.Item(1, 9).Value = 'Đơn nợ'
$ws3.Cells.Item(1, 10).Value = 'Nguồn khách'
$ws3.Cells.Item(1, 11).Value = 'Sale chính'
$ws3.Cells.Item(1, 12).Value = 'Đơn giá gốc'
$ws3.Cells.Item(1, 13).Value = 'Sale phụ'
$ws3.Cells.Item(1, 14).Value = 'Upsale'
$ws3.Cells.Item(1, 15).Value = 'Bác sĩ 1'
$ws3.Cells.Item(1, 16).Value = 'Bác sĩ 2'
$ws3.Cells.Item(1, 17).Value = 'Thanh toán lần đầu'
$ws3.Cells.Item(1, 18).Value = 'Đã thanh toán'
$ws3.Cells.Item(1, 19).Value = 'Tỉ lệ chiết khấu sale chính'
$ws3.Cells.Item(1, 20).Value = 'Tỉ lệ chiết khấu sale phụ'
$ws3.Cells.Item(1, 21).Value = 'id sale chính'
$ws3.Cells.Item(1, 22).Value = 'id sale phụ'
$ws3.Cells.Item(1, 23).Value = 'id bác sĩ 1'
$ws3.Cells.Item(1, 24).Value = 'id bác sĩ 2'
$ws3.Cells.Item(1, 25).Value = 'Chiết khấu bác sĩ 1'
$ws3.Cells.Item(1, 26).Value = 'Chiết khấu bác sĩ 2'
$ws3.Cells.Item(1, 27).Value = 'Chiết khấu sale chính'
$ws3.Cells.Item(1, 28).Value = 'Chiết khấu sale phụ'

# row 2 (data)
$ws3.Cells.Item(2, 1).Value = '06-30-2024'
$ws3.Cells.Item(2, 2).Value = '07-07-2024'
$ws3.Cells.Item(2, 3).Value = '59e6b0af-58c2-40e0-acb9-74a39fc79966'
$ws3.Cells.Item(2, 4).Value = 'TN'
$ws3.Cells.Item(2, 5).Value = 144
$ws3.Cells.Item(2, 6).Value = 'SÓC TRĂNG'
$ws3.Cells.Item(2, 7).Value = '7fd92044-b934-43de-9646-21004f1e6711'
$ws3.Cells.Item(2, 8).Value = 2000000
$ws3.Cells.Item(2, 9).Value = 'HD-LUXURY-502'
$ws3.Cells.Item(2, 10).Value = 'Cá nhân'
$ws3.Cells.Item(2, 11).Value = 'Lâm Thị Mỹ Hằng'
$ws3.Cells.Item(2, 12).Value = 8450000
$ws3.Cells.Item(2, 13).Value = 'Đỗ Thị Huyền Trân'
$ws3.Cells.Item(2, 14).Value = $null
$ws3.Cells.Item(2, 15).Value = 'Lâm Thị Mỹ Hằng'
$ws3.Cells.Item(2, 16).Value = $null
$ws3.Cells.Item(2, 17).Value = 2850000
$ws3.Cells.Item(2, 18).Value = 4850000
$ws3.Cells.Item(2, 19).Value = 0.13
$ws3.Cells.Item(2, 20).Value = 0.04
$ws3.Cells.Item(2, 21).Value = 'bc9b2b6b-3140-44b9-a1be-4dc8e77d8898'
$ws3.Cells.Item(2, 22).Value = 'f973382b-037a-4eb1-84bc-e9e5318184b8'
$ws3.Cells.Item(2, 23).Value = 'bc9b2b6b-3140-44b9-a1be-4dc8e77d8898'
$ws3.Cells.Item(2, 24).Value = $null
$ws3.Cells.Item(2, 25).Value = 200000
$ws3.Cells.Item(2, 26).Value = 0
$ws3.Cells.Item(2, 27).Value = 260000
$ws3.Cells.Item(2, 28).Value = 0

# row 3 (Tổng)
$ws3.Cells.Item(3, 4).Value = 'Tổng'
$ws3.Cells.Item(3, 5).Value = 1
$ws3.Cells.Item(3, 8).Value = 2000000
$ws3.Cells.Item(3, 12).Value = 8450000
$ws3.Cells.Item(3, 14).Value = 0
$ws3.Cells.Item(3, 17).Value = 2850000
$ws3.Cells.Item(3, 18).Value = 4850000
$ws3.Cells.Item(3, 19).Value = 0.13
$ws3.Cells.Item(3, 20).Value = 0.04
$ws3.Cells.Item(3, 25).Value = 200000
$ws3.Cells.Item(3, 26).Value = 0
$ws3.Cells.Item(3, 27).Value = 260000
$ws3.Cells.Item(3, 28).Value = 0

# --- Step 5: "Lương" sheet (new sheet) content ---
$ws4 = $sheetLuong

$ws4.Cells.Item(1, 1).Value = 'Danh mục'
$ws4.Cells.Item(1, 2).Value = 7
$ws4.Cells.Item(2, 1).Value = 'Ngày công'
$ws4.Cells.Item(2, 2).Value = 16
$ws4.Cells.Item(3, 1).Value = 'Phụ cấp'
$ws4.Cells.Item(3, 2).Value = 560000
$ws4.Cells.Item(4, 1).Value = 'Lương cơ bản tại CẦN THƠ'
$ws4.Cells.Item(4, 2).Value = 2571428.571428571
$ws4.Cells.Item(5, 1).Value = 'Chiết khấu sale chính tại CẦN THƠ'
$ws4.Cells.Item(5, 2).Value = 4236000
$ws4.Cells.Item(6, 1).Value = 'Chiết khấu sale phụ tại CẦN THƠ'
$ws4.Cells.Item(6, 2).Value = 248000
$ws4.Cells.Item(7, 1).Value = 'Đơn 1 bác sĩ tại CẦN THƠ'
$ws4.Cells.Item(7, 2).Value = 0
$ws4.Cells.Item(8, 1).Value = 'Đơn 2 bác sĩ tại CẦN THƠ'
$ws4.Cells.Item(8, 2).Value = 0
$ws4.Cells.Item(9, 1).Value = 'Công phụ phẫu 1 tại CẦN THƠ'
$ws4.Cells.Item(9, 2).Value = 0
$ws4.Cells.Item(10, 1).Value = 'Công phụ phẫu 2 tại CẦN THƠ'
$ws4.Cells.Item(10, 2).Value = 0
$ws4.Cells.Item(11, 1).Value = 'Ứng lương tại CẦN THƠ'
$ws4.Cells.Item(11, 2).Value = -350000
$ws4.Cells.Item(12, 1).Value = 'Lương cơ bản tại LONG XUYÊN'
$ws4.Cells.Item(12, 2).Value = $null
$ws4.Cells.Item(13, 1).Value = 'Chiết khấu sale chính tại LONG XUYÊN'
$ws4.Cells.Item(13, 2).Value = 0
$ws4.Cells.Item(14, 1).Value = 'Chiết khấu sale phụ tại LONG XUYÊN'
$ws4.Cells.Item(14, 2).Value = 0
$ws4.Cells.Item(15, 1).Value = 'Đơn 1 bác sĩ tại LONG XUYÊN'
$ws4.Cells.Item(15, 2).Value = 0
$ws4.Cells.Item(16, 1).Value = 'Đơn 2 bác sĩ tại LONG XUYÊN'
$ws4.Cells.Item(16, 2).Value = 0
$ws4.Cells.Item(17, 1).Value = 'Công phụ phẫu 1 tại LONG XUYÊN'
$ws4.Cells.Item(17, 2).Value = 0
$ws4.Cells.Item(18, 1).Value = 'Công phụ phẫu 2 tại LONG XUYÊN'
$ws4.Cells.Item(18, 2).Value = 0
$ws4.Cells.Item(19, 1).Value = 'Ứng lương tại LONG XUYÊN'
$ws4.Cells.Item(19, 2).Value = 0
$ws4.Cells.Item(20, 1).Value = 'Lương cơ bản tại SÓC TRĂNG'
$ws4.Cells.Item(20, 2).Value = $null
$ws4.Cells.Item(21, 1).Value = 'Chiết khấu sale chính tại SÓC TRĂNG'
$ws4.Cells.Item(21, 2).Value = 120000
$ws4.Cells.Item(22, 1).Value = 'Chiết khấu sale phụ tại SÓC TRĂNG'
$ws4.Cells.Item(22, 2).Value = 260000
$ws4.Cells.Item(23, 1).Value = 'Đơn 1 bác sĩ tại SÓC TRĂNG'
$ws4.Cells.Item(23, 2).Value = 0
$ws4.Cells.Item(24, 1).Value = 'Đơn 2 bác sĩ tại SÓC TRĂNG'
$ws4.Cells.Item(24, 2).Value = 0
$ws4.Cells.Item(25, 1).Value = 'Công phụ phẫu 1 tại SÓC TRĂNG'
$ws4.Cells.Item(25, 2).Value = 0
$ws4.Cells.Item(26, 1).Value = 'Công phụ phẫu 2 tại SÓC TRĂNG'
$ws4.Cells.Item(26, 2).Value = 0
$ws4.Cells.Item(27, 1).Value = 'Ứng lương tại SÓC TRĂNG'
$ws4.Cells.Item(27, 2).Value = 0
$ws4.Cells.Item(28, 1).Value = 'Tổng lương tại CẦN THƠ'
$ws4.Cells.Item(28, 2).Value = 11465428.57142857
$ws4.Cells.Item(29, 1).Value = 'Tổng lương tại LONG XUYÊN'
$ws4.Cells.Item(29, 2).Value = 0
$ws4.Cells.Item(30, 1).Value = 'Tổng lương tại SÓC TRĂNG'
$ws4.Cells.Item(30, 2).Value = 380000
$ws4.Cells.Item(31, 1).Value = 'Tổng lương'
$ws4.Cells.Item(31, 2).Value = 11845428.57142857

Write-Output "edit applied"
